# orders flow updated, Payment Methods & Import fields added
# Two new import fields ("units" and "expire_date") are inserted into the
# header row right after "quantity"; the remaining headers (price,
# description, status) shift two columns to the right so the full header
# row reads: branch_name, product_name, quantity, units, expire_date,
# price, description, status.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing trailing headers two slots to the right and insert
# the two new field names in their place.
$ws.Range("D1").Value = "units"
$ws.Range("E1").Value = "expire_date"
$ws.Range("F1").Value = "price"
$ws.Range("G1").Value = "description"
$ws.Range("H1").Value = "status"

# New header cells need the same bold / shaded header style as the rest
# of row 1 - copy formatting only (not values) from the already-styled F1.
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Match the column width that was applied to the newly visible column G
# (same width bucket as the other data columns).
$ws.Range("G1").EntireColumn.ColumnWidth = 11

# Move the active selection the same way the source workbook ended up.
$ws.Range("F2").Select() | Out-Null
